$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, shifting existing rows 60..101 down to 61..102.
$ws.Rows.Item(60).Insert()

# Populate the new row 60 with the added weekly price record.
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44942
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = 100112040
$ws.Range("G60").Value = "Cilantro"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 200
$ws.Range("K60").Value = 4500
$ws.Range("L60").Value = 5000
$ws.Range("M60").Value = 4750
$ws.Range("N60").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 2375
$ws.Range("Q60").Value = 2
$ws.Range("R60").Value = "Hortaliza"
